$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column (Price) values are plain text in the source data (not numbers),
# e.g. thousand-separated "51.027.30" or zero-padded "379.58". Force the
# Text number format before assigning so Excel does not auto-convert these
# numeric-looking strings into floating point numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.027.30'
$ws.Range("E2").Value = '  -1.47%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.937.55'
$ws.Range("E3").Value = '  -1.39%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '379.58'
$ws.Range("E5").Value = '  +0.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.11'
$ws.Range("E6").Value = '  -2.86%  '

# Row 7
$ws.Range("E7").Value = '  -1.61%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.586'
$ws.Range("E9").Value = '  -2.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.48'
$ws.Range("E10").Value = '  -2.97%  '

# Row 11
$ws.Range("E11").Value = '  -0.85%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0838'
$ws.Range("E12").Value = '  -0.86%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.414.10'
$ws.Range("E13").Value = '  -0.87%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.99'
$ws.Range("E14").Value = '  -3.87%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.37'
$ws.Range("E15").Value = '  -1.50%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.940.50'
$ws.Range("E16").Value = '  -1.52%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.984'
$ws.Range("E17").Value = '  +2.47%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '50.989.26'
$ws.Range("E18").Value = '  -1.66%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.19'
$ws.Range("E19").Value = '  -8.00%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.12'
$ws.Range("E20").Value = '  -4.14%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.52'
$ws.Range("E21").Value = '  -5.08%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0952'
$ws.Range("E22").Value = '  -0.97%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.47'
$ws.Range("E23").Value = '  -0.57%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '261.64'
$ws.Range("E24").Value = '  -0.84%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.91'
$ws.Range("E25").Value = '  +4.04%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.27'
$ws.Range("E26").Value = '  +10.71%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.68'
$ws.Range("E27").Value = '  +2.96%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.168'
$ws.Range("E28").Value = '  -1.17%  '

# Row 29
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.11%  '

# Row 30
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.112'
$ws.Range("E30").Value = '  +7.44%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.60'
$ws.Range("E31").Value = '  -1.89%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.78'
$ws.Range("E32").Value = '  -1.66%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.12'
$ws.Range("E33").Value = '  -1.59%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0455'
$ws.Range("E34").Value = '  +3.94%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.43'
$ws.Range("E35").Value = '  -2.50%  '

# Row 36
$ws.Range("E36").Value = '  -1.24%  '

# Row 37
$ws.Range("E37").Value = '  +0.07%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.97'
$ws.Range("E38").Value = '  -3.48%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.77'
$ws.Range("E39").Value = '  -4.04%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.57'
$ws.Range("E40").Value = '  -3.83%  '

# Row 41
$ws.Range("E41").Value = '  -0.83%  '

# Row 42
$ws.Range("E42").Value = '  -4.42%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.84'
$ws.Range("E43").Value = '  -2.59%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.23'
$ws.Range("E44").Value = '  -4.02%  '

# Row 45
$ws.Range("E45").Value = '  -1.60%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.39'
$ws.Range("E46").Value = '  +3.04%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.273'
$ws.Range("E47").Value = '  -3.49%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.23'
$ws.Range("E48").Value = '  -0.34%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.005.35'
$ws.Range("E49").Value = '  -1.61%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0346'
$ws.Range("E50").Value = '  +3.34%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.478'
$ws.Range("E51").Value = '  +12.95%  '
